# Error Calculations and Plots
# Apply edits to the "missing_data" worksheet:
#  1. Remove the two data rows whose ID is "RM 232" and "SC 92" (rows shift up).
#  2. Update a handful of B/C column values that became newly "missing" or
#     newly populated after the row removal / re-imputation pass.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Delete the "RM 232" row (originally row 26) ---------------------
$ws.Rows("26").Delete()

# After the above delete, the row that used to hold "SC 92" (originally
# row 28) has shifted up to row 27. Delete it too.
$ws.Rows("27").Delete()

# --- 2. Cell level value updates on the resulting (post-delete) sheet ---
# Row 6  (RM 21)  : C6  missing -> 15.1
$ws.Range("C6").Value = 15.1

# Row 8  (RM 38)  : C8  15.5 -> missing
$ws.Range("C8").ClearContents()

# Row 19 (RM 125) : C19 missing -> 13.2
$ws.Range("C19").Value = 13.2

# Row 21 (RM 135) : C21 12.7 -> missing
$ws.Range("C21").ClearContents()

# Row 23 (RM 140) : C23 missing -> 12.2
$ws.Range("C23").Value = 12.2

# Row 26 (SC 5)   : B26 -20.2 -> missing
$ws.Range("B26").ClearContents()

# Row 27 (SC 101) : B27 missing -> -20.4 ; C27 10 -> missing
$ws.Range("B27").Value = -20.4
$ws.Range("C27").ClearContents()

# Row 29 (SC 119) : B29 -19.5 -> missing ; C29 missing -> 11.2
$ws.Range("B29").ClearContents()
$ws.Range("C29").Value = 11.2
